# Applies the Alcam-Nrp1 NATMI re-run (Dr Hou advice): recompute the
# ECs/FAPs/sCs cross-tabulated ligand-receptor stats and extend the table
# from 3 to 9 sending/target-cluster combinations (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Alcam", "Nrp1", "ECs", 3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 2226.382456206273, 20037.44210585645, 0.406865147860101, 0.406865147860101),
  @("ECs", "Alcam", "Nrp1", "FAPs", 3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 1988.629009115025, 17897.66108203522, 0.3634163724103239, 0.3634163724103239),
  @("ECs", "Alcam", "Nrp1", "sCs", 3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 1135.049099478546, 10215.44189530691, 0.2074270386026733, 0.2074270386026733),
  @("FAPs", "Alcam", "Nrp1", "ECs", 2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 25.24666286081878, 227.219965747369, 0.004613756809485683, 0.004613756809485683),
  @("FAPs", "Alcam", "Nrp1", "FAPs", 2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 22.55059368098055, 202.955343128825, 0.004121057730565924, 0.004121057730565924),
  @("FAPs", "Alcam", "Nrp1", "sCs", 2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 12.87119464363756, 115.840751792738, 0.002352174711591665, 0.002352174711591665),
  @("sCs", "Alcam", "Nrp1", "ECs", 3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 25.51414208262167, 229.627278743595, 0.004662637886869745, 0.004662637886869745),
  @("sCs", "Alcam", "Nrp1", "FAPs", 3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 22.78950902920833, 205.105581262875, 0.004164718840188743, 0.004164718840188743),
  @("sCs", "Alcam", "Nrp1", "sCs", 3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 13.00756027524333, 117.06804247719, 0.002377095148200236, 0.002377095148200236)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $vals = $rows[$i]
  for ($c = 0; $c -lt $vals.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
  }
}